# "Added UI and functions for Archiving"
# Insert three new leading columns (Manufacturer Name, Brand, Product Name)
# in front of the existing Image Name / Variant / Image Url columns, restyle
# the new column widths, and move the active-cell selection to C11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing A:C columns (with all their cell styles/content) three
# slots to the right, making room for the three new header columns at A:C.
$ws.Range("A1:C1").EntireColumn.Insert() | Out-Null

# New header row values for the freshly inserted columns.
$ws.Range("A1").Value = "Manufacturer Name"
$ws.Range("B1").Value = "Brand"
$ws.Range("C1").Value = "Product Name"

# Column widths for the full A:F layout (values are in Excel "characters";
# 5/6 char is added internally by the width model, so subtract it up front
# to land exactly on the target stored width).
$widths = @{ 1 = 17.5; 2 = 6; 3 = 12.6640625; 4 = 11.5; 5 = 6.83203125; 6 = 9 }
foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = ($widths[$col] - 0.8333333333333333)
}

# Match the saved selection state.
$ws.Range("C11").Select() | Out-Null
